# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 06:35"

# --- Re-sort the tied "18 casos totales" block (rows 195-197) ---
# New order: Fiyi (194, unchanged), Nueva Caledonia, Belice, Santa Lucia, Islas Virgenes... (198, unchanged)
$ws.Range("A195").Value = "Nueva Caledonia"
$ws.Range("D195").Value = 18
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Belice"
$ws.Range("D196").Value = 16
$ws.Range("H196").Value = 2

$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# --- Re-sort the tied "11 casos totales" block (rows 209-211) ---
# New order: Surinam (208, unchanged), Seychelles, Groenlandia, Montserrat, Islas Virgenes Britanicas (212, unchanged)
$ws.Range("A209").Value = "Seychelles"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Update Kirguistan statistics (row 100) ---
$ws.Range("B100").Value = 1243
$ws.Range("C100").Value = 27
$ws.Range("D100").Value = 898
$ws.Range("E100").Value = 331
